$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove old rows 8:10 (previously the MuSCs-sending-cluster block), shrinking the table to 6 data rows
$ws.Rows("8:10").Delete()

# Rewrite rows 2:7 with the updated TPM-based values
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Has2"
$ws.Range("C2").Value = "Cd44"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 17.78255766666667
$ws.Range("H2").Value = 53.347673
$ws.Range("I2").Value = 0.9887364018488972
$ws.Range("J2").Value = 0.9887364018488973
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 7.487621999999999
$ws.Range("N2").Value = 22.462866
$ws.Range("O2").Value = 0.1384395179233961
$ws.Range("P2").Value = 0.1384395179233961
$ws.Range("Q2").Value = 133.149070001202
$ws.Range("R2").Value = 1198.341630010818
$ws.Range("S2").Value = 0.1368801908252745
$ws.Range("T2").Value = 0.1368801908252746

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Has2"
$ws.Range("C3").Value = "Cd44"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 17.78255766666667
$ws.Range("H3").Value = 53.347673
$ws.Range("I3").Value = 0.9887364018488972
$ws.Range("J3").Value = 0.9887364018488973
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 31.999428
$ws.Range("N3").Value = 95.998284
$ws.Range("O3").Value = 0.5916411627275552
$ws.Range("P3").Value = 0.5916411627275552
$ws.Range("Q3").Value = 569.031673710348
$ws.Range("R3").Value = 5121.285063393132
$ws.Range("S3").Value = 0.5849771544209408
$ws.Range("T3").Value = 0.5849771544209409

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Has2"
$ws.Range("C4").Value = "Cd44"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 17.78255766666667
$ws.Range("H4").Value = 53.347673
$ws.Range("I4").Value = 0.9887364018488972
$ws.Range("J4").Value = 0.9887364018488973
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 14.59882166666667
$ws.Range("N4").Value = 43.796465
$ws.Range("O4").Value = 0.2699193193490487
$ws.Range("P4").Value = 0.2699193193490487
$ws.Range("Q4").Value = 259.6043881528827
$ws.Range("R4").Value = 2336.439493375945
$ws.Range("S4").Value = 0.2668790566026819
$ws.Range("T4").Value = 0.2668790566026819

$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Has2"
$ws.Range("C5").Value = "Cd44"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.2025773333333333
$ws.Range("H5").Value = 0.6077319999999999
$ws.Range("I5").Value = 0.01126359815110275
$ws.Range("J5").Value = 0.01126359815110275
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 7.487621999999999
$ws.Range("N5").Value = 22.462866
$ws.Range("O5").Value = 0.1384395179233961
$ws.Range("P5").Value = 0.1384395179233961
$ws.Range("Q5").Value = 1.516822497768
$ws.Range("R5").Value = 13.651402479912
$ws.Range("S5").Value = 0.00155932709812152
$ws.Range("T5").Value = 0.00155932709812152

$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Has2"
$ws.Range("C6").Value = "Cd44"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.2025773333333333
$ws.Range("H6").Value = 0.6077319999999999
$ws.Range("I6").Value = 0.01126359815110275
$ws.Range("J6").Value = 0.01126359815110275
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 31.999428
$ws.Range("N6").Value = 95.998284
$ws.Range("O6").Value = 0.5916411627275552
$ws.Range("P6").Value = 0.5916411627275552
$ws.Range("Q6").Value = 6.482358792431999
$ws.Range("R6").Value = 58.34122913188799
$ws.Range("S6").Value = 0.006664008306614371
$ws.Range("T6").Value = 0.006664008306614371

$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Has2"
$ws.Range("C7").Value = "Cd44"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.2025773333333333
$ws.Range("H7").Value = 0.6077319999999999
$ws.Range("I7").Value = 0.01126359815110275
$ws.Range("J7").Value = 0.01126359815110275
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 14.59882166666667
$ws.Range("N7").Value = 43.796465
$ws.Range("O7").Value = 0.2699193193490487
$ws.Range("P7").Value = 0.2699193193490487
$ws.Range("Q7").Value = 2.957390363042222
$ws.Range("R7").Value = 26.61651326738
$ws.Range("S7").Value = 0.003040262746366857
$ws.Range("T7").Value = 0.003040262746366857

